$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '26.033.34'
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.629.87'
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -0.81%  '
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '214.27'
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -0.49%  '
$c.Style = 'Normal'

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -0.65%  '
$c.Style = 'Normal'

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.12%  '
$c.Style = 'Normal'

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -1.65%  '
$c.Style = 'Normal'

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -2.95%  '
$c.Style = 'Normal'

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '18.51'
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -5.52%  '
$c.Style = 'Normal'

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0789'
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  -0.98%  '
$c.Style = 'Normal'

$c = $ws.Range('B12')
$c.NumberFormat = '@'
$c.Value = 'WrappedliquidstakedEther2.0'
$c.Style = 'Normal'
$c = $ws.Range('C12')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.857.88'
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -0.72%  '
$c.Style = 'Normal'

$c = $ws.Range('B13')
$c.NumberFormat = '@'
$c.Value = 'WrappedEther'
$c.Style = 'Normal'
$c = $ws.Range('C13')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '1.628.55'
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -0.30%  '
$c.Style = 'Normal'

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '4.19'
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -1.76%  '
$c.Style = 'Normal'

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  -2.96%  '
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '26.040.33'
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -0.05%  '
$c.Style = 'Normal'

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.0₃0742'
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  -2.38%  '
$c.Style = 'Normal'

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '61.55'
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -2.98%  '
$c.Style = 'Normal'

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '192.69'
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -1.01%  '
$c.Style = 'Normal'

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -2.34%  '
$c.Style = 'Normal'

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -3.50%  '
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '6.06'
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  -2.19%  '
$c.Style = 'Normal'

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +1.28%  '
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '144.33'
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +0.24%  '
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +0.33%  '
$c.Style = 'Normal'

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -3.50%  '
$c.Style = 'Normal'

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -2.30%  '
$c.Style = 'Normal'

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '15.27'
$c.Style = 'Normal'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -1.54%  '
$c.Style = 'Normal'

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -0.59%  '
$c.Style = 'Normal'

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -2.51%  '
$c.Style = 'Normal'

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -3.98%  '
$c.Style = 'Normal'

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  -4.74%  '
$c.Style = 'Normal'

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  -2.83%  '
$c.Style = 'Normal'

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.42'
$c.Style = 'Normal'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -1.51%  '
$c.Style = 'Normal'

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.127.43'
$c.Style = 'Normal'
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -0.27%  '
$c.Style = 'Normal'

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -5.73%  '
$c.Style = 'Normal'

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -0.88%  '
$c.Style = 'Normal'

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -3.35%  '
$c.Style = 'Normal'

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  -2.02%  '
$c.Style = 'Normal'

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -0.47%  '
$c.Style = 'Normal'

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '1.766.98'
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -0.80%  '
$c.Style = 'Normal'

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -4.86%  '
$c.Style = 'Normal'

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -5.73%  '
$c.Style = 'Normal'

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -1.45%  '
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '54.49'
$c.Style = 'Normal'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -3.53%  '
$c.Style = 'Normal'

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.0523'
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +0.38%  '
$c.Style = 'Normal'

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -0.55%  '
$c.Style = 'Normal'

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '7.53'
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -3.27%  '
$c.Style = 'Normal'

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  +0.24%  '
$c.Style = 'Normal'
